$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Domino LED Rev. B"
